# Insert a new data row at row 268 (pushing the existing rows 268..358 down
# to 269..359), then populate the new row with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(268).Insert()

$ws.Cells.Item(268, 1).Value2 = 4
$ws.Cells.Item(268, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(268, 3).Value2 = "Los Lagos"
$ws.Cells.Item(268, 4).Value2 = 44988
$ws.Cells.Item(268, 5).Value2 = 10
$ws.Cells.Item(268, 6).Value2 = 100112032
$ws.Cells.Item(268, 7).Value2 = "Zapallo italiano"
$ws.Cells.Item(268, 8).Value2 = "Sin especificar"
$ws.Cells.Item(268, 9).Value2 = "Primera"
$ws.Cells.Item(268, 10).Value2 = 250
$ws.Cells.Item(268, 11).Value2 = 15000
$ws.Cells.Item(268, 12).Value2 = 15000
$ws.Cells.Item(268, 13).Value2 = 15000
$ws.Cells.Item(268, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(268, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(268, 16).Value2 = 300
$ws.Cells.Item(268, 17).Value2 = 50
$ws.Cells.Item(268, 18).Value2 = "Hortaliza"

# Keep the date-formatted style (same as D269 and the other D-column cells).
$ws.Cells.Item(268, 4).NumberFormat = $ws.Cells.Item(269, 4).NumberFormat
